# Applies the "cryptos list" refresh: updates Price (D) and Volume(1h) (E)
# columns for rows 2-51 to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.448.37"
$ws.Range("E2").Value = "  -0.78%  "

$ws.Range("D3").Value = "2.302.53"
$ws.Range("E3").Value = "  -1.86%  "

$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.42"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.28"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.49%  "

$ws.Range("E7").Value = "  -0.56%  "

$ws.Range("E8").Value = "  -0.21%  "

$ws.Range("E9").Value = "  -0.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.70"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0909"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.97%  "

$ws.Range("E12").Value = "  -1.92%  "

$ws.Range("E13").Value = "  -0.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.966"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.16%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.32"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.25%  "

$ws.Range("D16").Value = "2.649.28"
$ws.Range("E16").Value = "  -2.05%  "

$ws.Range("D17").Value = "2.300.81"
$ws.Range("E17").Value = "  -2.26%  "

$ws.Range("D18").Value = "42.239.77"
$ws.Range("E18").Value = "  -1.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.44"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.16%  "

$ws.Range("E20").Value = "  -0.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.66"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.61"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.86%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "280.53"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.33"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +19.31%  "

$ws.Range("E25").Value = "  -2.41%  "

$ws.Range("E26").Value = "  -0.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.97"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.19%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.89"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.39"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.88%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "23.02"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.31"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "164.49"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0875"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.91"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.98%  "

$ws.Range("E35").Value = "  +3.74%  "

$ws.Range("E36").Value = "  -9.60%  "

$ws.Range("E37").Value = "  -5.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.61"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.62%  "

$ws.Range("E39").Value = "  -2.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.78"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.46%  "

$ws.Range("E41").Value = "  +4.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.05"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.43%  "

$ws.Range("E43").Value = "  -1.86%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "69.24"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.07%  "

$ws.Range("E45").Value = "  -4.52%  "

$ws.Range("E46").Value = "  -0.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.02"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.63%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "112.13"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.90%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "77.51"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.92"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.61%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.28"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.13%  "
